$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.792.06'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.739.71'
$ws.Range('E3').Value = '  +19.16%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '614.92'
$ws.Range('E5').Value = '  +6.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.21'
$ws.Range('E6').Value = '  -2.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.737.10'
$ws.Range('E7').Value = '  +19.10%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  +4.14%  '
$ws.Range('E10').Value = '  +9.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.42'
$ws.Range('E11').Value = '  -1.86%  '
$ws.Range('E12').Value = '  +7.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '41.00'
$ws.Range('E13').Value = '  +11.23%  '
$ws.Range('E14').Value = '  +5.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.357.45'
$ws.Range('E15').Value = '  +19.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.739.23'
$ws.Range('E16').Value = '  +19.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.824.96'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('E19').Value = '  +6.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '517.33'
$ws.Range('E20').Value = '  +6.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.73'
$ws.Range('E21').Value = '  +0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.39'
$ws.Range('E22').Value = '  +20.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.736'
$ws.Range('E23').Value = '  +5.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.83'
$ws.Range('E24').Value = '  +5.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.48'
$ws.Range('E25').Value = '  +6.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.55'
$ws.Range('E26').Value = '  +4.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.90'
$ws.Range('E27').Value = '  +2.83%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('E29').Value = '  +31.62%  '
$ws.Range('E30').Value = '  +6.12%  '
$ws.Range('E31').Value = '  +8.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.86'
$ws.Range('E32').Value = '  -3.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.58'
$ws.Range('E33').Value = '  +11.67%  '
$ws.Range('E34').Value = '  +2.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.23'
$ws.Range('E36').Value = '  +9.73%  '
$ws.Range('E37').Value = '  +8.11%  '
$ws.Range('E38').Value = '  +5.31%  '
$ws.Range('E39').Value = '  +6.47%  '
$ws.Range('E40').Value = '  +6.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.44'
$ws.Range('E41').Value = '  +4.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.88'
$ws.Range('E42').Value = '  +5.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '44.60'
$ws.Range('E43').Value = '  -8.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '421.83'
$ws.Range('E44').Value = '  +5.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.076.72'
$ws.Range('E45').Value = '  +9.93%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('E47').Value = '  +4.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.90'
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '136.61'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.51'
$ws.Range('E50').Value = '  +6.41%  '
$ws.Range('E51').Value = '  -0.02%  '
